# Roadmap workbook update: "Update de textos webinars y fotos"
# - Adds a new "Sheet2" with a priority/improvement table.
# - Makes Sheet2 the active/selected sheet, Sheet1 keeps a specific selection.
# - Increases the height of a few rows on Sheet1 (rows 3, 4, 7).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: selection + row-height tweaks -----------------------------
$ws1.Select()
$ws1.Range("D3:E8").Select() | Out-Null

$ws1.Rows.Item(3).RowHeight = 30
$ws1.Rows.Item(4).RowHeight = 30
$ws1.Rows.Item(7).RowHeight = 30

# --- Create Sheet2 right after Sheet1 -----------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "#"
$ws2.Range("B1").Value = "Mejora"
$ws2.Range("C1").Value = "Impacto(ventas)"
$ws2.Range("D1").Value = "Esfuerzo(técnico)"
$ws2.Range("E1").Value = "Urgencia"
$ws2.Range("F1").Value = "Prioridad sugerida"

# Data rows
$data = @(
    @(1,  "Mejorar descripción y estructura de clases individuales", 5, 2, 5, "1️⃣ inmediata"),
    @(2,  "Página dedicada de venta de módulos", 5, 3, 5, "2️⃣ muy alta"),
    @(7,  "Upsell en correo de confirmación", 4, 2, 4, "3️⃣ alta"),
    @(5,  "Opción de asesoría personalizada desde páginas de webinars", 4, 3, 4, "4️⃣ alta"),
    @(6,  "En prelobby, enlace a asesoría 1-a-1", 3, 2, 4, "5️⃣ media"),
    @(4,  "Página de venta de uno-a-uno", 4, 4, 3, "6️⃣ media"),
    @(9,  "Duraciones 30/60/120 min en asesorías", 3, 2, 3, "7️⃣ media-baja"),
    @(8,  "Servicios adicionales (personalizar/configurar plantillas Excel)", 3, 3, 3, "8️⃣ media-baja"),
    @(10, "Selector de fechas en clases", 2, 4, 2, "9️⃣ baja"),
    @(3,  "Destacar próximo módulo en home", 2, 3, 2, "🔟 baja")
)

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Column widths (best-fit-like)
$ws2.Range("A1:F11").EntireColumn.AutoFit() | Out-Null

# Leave Sheet2 selected/active, matching the target selection
$ws2.Select()
$ws2.Range("B15").Select() | Out-Null
